$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-coerced to a number by Excel
# (losing fixed decimal formatting like trailing zeros) are first forced to
# Text format, written as text, then restored to the default "Normal" style
# so no stray number-format style is left attached to the cell.
$forceTextCells = @{
    'D5' = '307.59'
    'D6' = '96.49'
    'D9' = '0.502'
    'D10' = '35.33'
    'D12' = '18.40'
    'D17' = '0.782'
    'D19' = '12.95'
    'D22' = '67.19'
    'D23' = '235.92'
    'D26' = '1.00'
    'D28' = '25.11'
    'D30' = '166.39'
    'D32' = '33.00'
    'D34' = '4.77'
    'D35' = '4.99'
    'D36' = '17.63'
    'D38' = '0.0695'
    'D44' = '0.0280'
    'D45' = '18.27'
    'D46' = '10.08'
    'D47' = '2.05'
    'D48' = '2.81'
    'D49' = '2.92'
    'D50' = '54.02'
}
foreach ($ref in $forceTextCells.Keys) {
    $ws.Range($ref).NumberFormat = "@"
}
foreach ($ref in $forceTextCells.Keys) {
    $ws.Range($ref).Value = $forceTextCells[$ref]
}
foreach ($ref in $forceTextCells.Keys) {
    $ws.Range($ref).Style = "Normal"
}

# Remaining cells are naturally text (contain letters, "%", multiple "."
# separators, or surrounding spaces) so a direct .Value assignment keeps
# them as text without any extra formatting step.
$ws.Range('D2').Value = '42.848.71'
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').Value = '2.298.32'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  +2.45%  '
$ws.Range('E6').Value = '  -1.78%  '
$ws.Range('E7').Value = '  -2.31%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -2.46%  '
$ws.Range('E10').Value = '  -2.42%  '
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('E12').Value = '  +3.91%  '
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('D15').Value = '2.656.65'
$ws.Range('D16').Value = '2.295.49'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').Value = '42.750.03'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('E19').Value = '  +1.30%  '
$ws.Range('E20').Value = '  -1.71%  '
$ws.Range('E21').Value = '  -1.09%  '
$ws.Range('E22').Value = '  -2.55%  '
$ws.Range('E23').Value = '  -0.64%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('E28').Value = '  +0.79%  '
$ws.Range('E29').Value = '  +16.89%  '
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('E35').Value = '  -1.74%  '
$ws.Range('E36').Value = '  -1.59%  '
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('E40').Value = '  -1.42%  '
$ws.Range('E41').Value = '  -1.26%  '
$ws.Range('E42').Value = '  -2.87%  '
$ws.Range('D43').Value = '2.011.54'
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('E44').Value = '  -2.23%  '
$ws.Range('E45').Value = '  +4.51%  '
$ws.Range('E46').Value = '  -2.50%  '
$ws.Range('E47').Value = '  -8.09%  '
$ws.Range('E48').Value = '  -0.67%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E49').Value = '  +9.01%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('D51').Value = '2.522.26'
$ws.Range('E51').Value = '  -0.08%  '
